# Auto-generated edit script: apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.167.24'
$ws.Range("E2").Value = '  +2.47%  '
$ws.Range("D3").Value = '2.360.88'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.85'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.70'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.83'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0916'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.975'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -3.99%  '
$ws.Range("D15").Value = '2.716.86'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.22'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").Value = '2.360.39'
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '45.124.13'
$ws.Range("E18").Value = '  +2.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.29'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +10.18%  '
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.18'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -5.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.25'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.51'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.14'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -3.33%  '
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E27").Value = '  -1.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.19'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -4.92%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0976'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +7.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.33'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.23'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -5.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '167.74'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.00'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +4.92%  '
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.67'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.95'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +4.52%  '
$ws.Range("E39").Value = '  -3.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.87'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.78'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +2.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.79'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -5.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.66'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("E44").Value = '  -4.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.83'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -6.74%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.843.60'
$ws.Range("E46").Value = '  +11.04%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '83.53'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +7.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.70'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +7.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.64'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -4.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.19'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +1.93%  '
